# Update cryptos list values (prices and percent-volume changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.882.73"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.556.35"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.81"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.04"
$ws.Range("E6").Value = "  +6.52%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.41"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  +6.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.64"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.512.77"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.10"
$ws.Range("E15").Value = "  +6.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.886"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.937.46"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.70"
$ws.Range("E18").Value = "  +7.44%  "
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.78"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.47"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.98"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "27.90"
$ws.Range("E25").Value = "  -6.13%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.16"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.34"
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.02"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.54"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0810"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.34"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.55"
$ws.Range("E36").Value = "  +11.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.60"
$ws.Range("E37").Value = "  +11.94%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.12"
$ws.Range("E40").Value = "  +33.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0306"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.084.47"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.71"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.07"
$ws.Range("E47").Value = "  +3.87%  "

# Row 48: was RocketPoolETH -> now ordi
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "75.44"
$ws.Range("E48").Value = "  +8.35%  "

# Row 49: was ordi -> now RocketPoolETH
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.805.50"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.10"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.192"
$ws.Range("E51").Value = "  +2.15%  "
